$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -6
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = 9
$ws.Range("F10").Value = -8
$ws.Range("F12").Value = -6
$ws.Range("F16").Value = 3
$ws.Range("F18").Value = -8
